$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")

# --- Row 6 : "Scrapping python" / "Web Wikipédia" -> "Open Data " / "Site gouv",
#     candidate list URL moved to D6, and a new comment added in E6 ---
$ws.Range("A6").Value = "Open Data "
$ws.Range("C6").Value = "Site gouv"
$ws.Range("D6").Value = "Obtenir la liste des candidats des 577 circonsciptions et leurs infos "
$ws.Range("E6").Value = "Infos : Sexe, DatedeNaissance, Nuance, Profession, FonctionPublique, Sortant + Infos sur le suppléant"

# --- Row 12 : new comment added in E12 (rest unchanged) ---
$ws.Range("E12").Value = "Infos : Sexe, (Inscrits, Votants, Blancs) par circo, NbrVoix"

# --- Row 14 : D14 text replaced, new comment added in E14 ---
$ws.Range("D14").Value = "Récupérer les résultats des élections législatives 2007 et 2012"
$ws.Range("E14").Value = "Infos :  Sexe, (Inscrits, Votants, Blancs) par circo, Nuance, NbrVoix"

# --- Row 16 : D16 text replaced ---
$ws.Range("D16").Value = "Obtenir la liste des partis et leurs infos"

# --- Row 18 : previously empty styled row, now filled with new data ---
$ws.Range("A18").Value = "Open Data "
$ws.Range("B18").Value = "Données froides"
$ws.Range("C18").Value = "INSEE"
$ws.Range("D18").Value = "Obtenir des données extérieures qui ont une influence sur les élections"
$ws.Range("E18").Value = "Infos : Taux du chômage trimestrielle par département (2016-1982)"

# --- Update the view: scroll so row 4 is the top-left visible row, and
#     move the active selection to D29 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D29").Select()
